# Update countries & provincias Spain
# Refresh the COVID stats dump on sheet "Pais": update the "last updated"
# timestamp and the per-country rows whose figures changed between the two
# data pulls. Because the sheet is kept sorted by "Casos totales"
# descending, a handful of rows also swap which country occupies them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner in A1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 5 de Septiembre de 2020 a las 00:06"

# row => Pais, Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes
$rows = @(
    @(4,   "Estados Unidos",      6381221, 45977, 3620974, 2568417, 0, 772, 191830),
    @(5,   "Brasil",              4091801, 45651, 3247610,  718689, 0, 773, 125502),
    @(12,  "España",               517133,  4503,       0,       0, 0, 184,  29418),
    @(19,  "Francia",              309156,  8975,   87206,  191226, 0,  18,  30724),
    @(81,  "Costa de Marfil",       18269,    61,   17162,     988, 0,   0,    119),
    @(95,  "Guinea",                 9649,    70,    8832,     756, 0,   0,     61),
    @(105, "Mauritania",             7126,    20,    6623,     343, 0,   0,    160),
    @(108, "Malaui",                 5608,    15,    3533,    1900, 0,   0,    175),
    @(111, "Guinea Ecuatorial",      4972,     7,    4413,     476, 0,   0,     83),
    @(118, "Ruanda",                 4304,    49,    2191,    2095, 0,   0,     18),
    @(119, "Eslovaquia",             4300,   137,    2693,    1570, 0,   0,     37),
    @(120, "Cuba",                   4266,    52,    3487,     679, 0,   0,    100),
    @(121, "Mozambique",             4265,    58,    2511,    1728, 0,   0,     26),
    @(123, "Cabo Verde",             4200,    75,    3562,     597, 0,   0,     41),
    @(131, "Siria",                  3041,    68,     698,    2216, 0,   3,    127),
    @(132, "Lituania",               3004,    26,    1920,     998, 0,   0,     86),
    @(135, "Mali",                   2814,     7,    2220,     468, 0,   0,    126),
    @(144, "Sierra Leona",           2041,     6,    1602,     368, 0,   0,     71),
    @(154, "Togo",                   1457,    14,    1079,     347, 0,   0,     31),
    @(214, "Montserrat",               13,     0,      12,       0, 0,   0,      1),
    @(215, "Islas Malvinas",           13,     0,      13,       0, 0,   0,      0)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
    $ws.Cells.Item($rowNum, 7).Value = $r[7]
    $ws.Cells.Item($rowNum, 8).Value = $r[8]
}
